$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 5000
$ws.Range("I6").Value = 5000
$ws.Range("K6").Value = 15000
$ws.Range("M6").Value = -14888

$ws.Range("H18").Value = 857.1429000000001
$ws.Range("I18").Value = 857.1429000000001
$ws.Range("K18").Value = 857.1429000000001
$ws.Range("M18").Value = -573.1429000000001

$ws.Range("H21").Value = 717
$ws.Range("I21").Value = 717
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 717
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -249
$ws.Range("N21").Value = ""

$ws.Range("H23").Value = 717
$ws.Range("I23").Value = 717
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 717
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -483
$ws.Range("N23").Value = ""

$ws.Range("H136").Value = 200000
$ws.Range("J136").Value = 200000
$ws.Range("L136").Value = 200000
$ws.Range("N136").Value = -210200

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 154.28572
$ws.Range("I4").Value = 159.5
$ws.Range("K4").Value = 159.5
$ws.Range("M4").Value = -43.5

$ws.Range("H6").Value = 3006.111
$ws.Range("I6").Value = 5250.5
$ws.Range("J6").Value = 2364.8572
$ws.Range("K6").Value = 5250.5
$ws.Range("L6").Value = 2364.8572
$ws.Range("M6").Value = -5077.5
$ws.Range("N6").Value = -2710.8572

$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = ""

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = ""

$ws.Range("H46").Value = 19500
$ws.Range("I46").Value = 19000
$ws.Range("K46").Value = 19000
$ws.Range("M46").Value = -18681

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 33681.25
$ws.Range("I26").Value = 34821.43
$ws.Range("K26").Value = 34821.43
$ws.Range("M26").Value = -34529.43

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 639.8570999999999
$ws.Range("I16").Value = 639.8570999999999
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 639.8570999999999
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -352.8570999999999
$ws.Range("N16").Value = ""

$ws.Range("H58").Value = 805.64703
$ws.Range("I58").Value = 833.13336
$ws.Range("K58").Value = 833.13336
$ws.Range("M58").Value = -630.13336

$ws.Range("H86").Value = 5551
$ws.Range("I86").Value = 1502.3334
$ws.Range("K86").Value = 1502.3334
$ws.Range("M86").Value = -379.3334

$ws.Range("H89").Value = 5551
$ws.Range("I89").Value = 1502.3334
$ws.Range("K89").Value = 7511.666999999999
$ws.Range("M89").Value = -1895.666999999999

$ws.Range("H113").Value = 639.8570999999999
$ws.Range("I113").Value = 639.8570999999999
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 639.8570999999999
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1530.1429
$ws.Range("N113").Value = ""

$ws.Range("H136").Value = 805.64703
$ws.Range("I136").Value = 833.13336
$ws.Range("K136").Value = 2499.40008
$ws.Range("M136").Value = 50.59991999999966

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 79.57143000000001
$ws.Range("I7").Value = 90.40000000000001
$ws.Range("J7").Value = 52.5
$ws.Range("K7").Value = 271.2
$ws.Range("L7").Value = 157.5
$ws.Range("M7").Value = -159.2
$ws.Range("N7").Value = -381.5

$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = ""

$ws.Range("H46").Value = 3888.4443
$ws.Range("J46").Value = 5332.1665
$ws.Range("L46").Value = 15996.4995
$ws.Range("N46").Value = -16178.4995

$ws.Range("H55").Value = 4200
$ws.Range("I55").Value = 4200
$ws.Range("K55").Value = 12600
$ws.Range("M55").Value = -12423

$ws.Range("H68").Value = 1176.2
$ws.Range("I68").Value = 865.6667
$ws.Range("J68").Value = 1642
$ws.Range("K68").Value = 2597.0001
$ws.Range("L68").Value = 4926
$ws.Range("M68").Value = -1786.0001
$ws.Range("N68").Value = -6548

$ws.Range("H71").Value = 1176.2
$ws.Range("I71").Value = 865.6667
$ws.Range("J71").Value = 1642
$ws.Range("K71").Value = 7791.0003
$ws.Range("L71").Value = 14778
$ws.Range("M71").Value = -3735.0003
$ws.Range("N71").Value = -22890

$ws.Range("H80").Value = 75
$ws.Range("J80").Value = 100
$ws.Range("L80").Value = 300
$ws.Range("N80").Value = -2172

$ws.Range("H83").Value = 75
$ws.Range("J83").Value = 100
$ws.Range("L83").Value = 900
$ws.Range("N83").Value = -10260

$ws.Range("H92").Value = 524.2857
$ws.Range("I92").Value = 482.5
$ws.Range("K92").Value = 1447.5
$ws.Range("M92").Value = -199.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = ""

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1433.6666
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1433.6666
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1433.6666
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = -2023.6666

$ws.Range("H27").Value = 1433.6666
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1433.6666
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1433.6666
$ws.Range("M27").Value = ""
$ws.Range("N27").Value = -1647.6666

$ws.Range("H68").Value = 4766.3335
$ws.Range("I68").Value = 3999.5
$ws.Range("J68").Value = 5149.75
$ws.Range("K68").Value = 3999.5
$ws.Range("L68").Value = 5149.75
$ws.Range("M68").Value = -3250.5
$ws.Range("N68").Value = -6647.75

$ws.Range("H71").Value = 4766.3335
$ws.Range("I71").Value = 3999.5
$ws.Range("J71").Value = 5149.75
$ws.Range("K71").Value = 19997.5
$ws.Range("L71").Value = 25748.75
$ws.Range("M71").Value = -16253.5
$ws.Range("N71").Value = -33236.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 60000
$ws.Range("I40").Value = 60000
$ws.Range("K40").Value = 60000
$ws.Range("M40").Value = -59851

$ws.Range("H62").Value = 2400.3333
$ws.Range("J62").Value = 2000
$ws.Range("L62").Value = 2000
$ws.Range("N62").Value = -3248

$ws.Range("H65").Value = 2400.3333
$ws.Range("J65").Value = 2000
$ws.Range("L65").Value = 10000
$ws.Range("N65").Value = -16240

$ws.Range("H70").Value = 90000
$ws.Range("I70").Value = 90000
$ws.Range("K70").Value = 90000
$ws.Range("M70").Value = -89685

$ws.Range("H73").Value = 90000
$ws.Range("I73").Value = 90000
$ws.Range("K73").Value = 90000
$ws.Range("M73").Value = -88908

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = ""

$ws.Range("H140").Value = 157750
$ws.Range("J140").Value = 157750
$ws.Range("L140").Value = 157750
$ws.Range("N140").Value = -168110
